$wb = $excel.ActiveWorkbook

# Insert a new worksheet "BirdTextField" right after "BirdSuborderTranslation"
# (i.e. before "BirdFamily") to back a new quiz-game text-field table.
$afterSheet = $wb.Worksheets.Item("BirdSuborderTranslation")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "BirdTextField"

# Header row
$newSheet.Range("A1").Value = "bird"
$newSheet.Range("B1").Value = "language"
$newSheet.Range("C1").Value = "description"

# Sample data row
$newSheet.Range("A2").Value = "parus major"
$newSheet.Range("B2").Value = "eng"
$newSheet.Range("C2").Value = "Test description"

# Match the author's final selection/active cell on the new sheet
$newSheet.Range("C2").Select()
